$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the two "special" fonts currently used for the stimulus
#     cells (B2/B3 = hiragana tfb 60pt, C2/C3 = bold Calibri 60pt) before
#     we overwrite anything, so we can move them onto column A instead. ---
$hiraganaFont = $ws.Range("B2").Font
$hiraganaName = $hiraganaFont.Name
$hiraganaSize = $hiraganaFont.Size
$hiraganaBold = $hiraganaFont.Bold

$boldFont = $ws.Range("C2").Font
$boldName = $boldFont.Name
$boldSize = $boldFont.Size
$boldBold = $boldFont.Bold

# --- Row 1 header labels stay the same text, just rewritten. ---
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "learnt_letters"
$ws.Range("C1").Value = "numbers"

# --- Row 2 / Row 3 label cells keep their text but now pick up the big
#     custom fonts that used to live on the letter-stimulus columns. ---
$ws.Range("A2").Value = "Most Probable"
$ws.Range("A2").Font.Name = $hiraganaName
$ws.Range("A2").Font.Size = $hiraganaSize
$ws.Range("A2").Font.Bold = $hiraganaBold

$ws.Range("A3").Value = "Least Probable"
$ws.Range("A3").Font.Name = $boldName
$ws.Range("A3").Font.Size = $boldSize
$ws.Range("A3").Font.Bold = $boldBold

# --- The letter-stimulus columns (B/C) now hold picture filenames
#     instead of letters, rendered in plain default formatting. ---
$ws.Range("B2").Value = "Hselect3.jpg"
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Bold = $false

$ws.Range("C2").Value = "Nselect3.jpg"
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Size = 11
$ws.Range("C2").Font.Bold = $false

$ws.Range("B3").Value = "Hselect3.jpg"
$ws.Range("B3").Font.Name = "Calibri"
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Bold = $false

$ws.Range("C3").Value = "Nselect3.jpg"
$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").Font.Bold = $false

# --- Rows shrink now that they no longer host giant 60pt glyphs. ---
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 37.5

# --- View scrolled over to column B, selection parked at B11. ---
$ws.Range("B11").Select()
$excel.ActiveWindow.ScrollColumn = 2
